$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.848.62'
$ws.Range("E2").Value = '  -0.70%  '
$ws.Range("D3").Value = '1.839.95'
$ws.Range("E3").Value = '  +1.67%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = "'231.41"
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = "'39.58"
$ws.Range("E8").Value = '  -1.94%  '
$ws.Range("E9").Value = '  +1.97%  '
$ws.Range("D10").Value = "'0.0686"
$ws.Range("E10").Value = '  +0.26%  '
$ws.Range("D11").Value = "'0.0980"
$ws.Range("E11").Value = '  -1.69%  '
$ws.Range("D12").Value = '2.105.44'
$ws.Range("D13").Value = "'11.39"
$ws.Range("E13").Value = '  +3.30%  '
$ws.Range("D14").Value = '1.842.71'
$ws.Range("E14").Value = '  +1.95%  '
$ws.Range("D15").Value = "'0.672"
$ws.Range("E15").Value = '  +1.29%  '
$ws.Range("D16").Value = "'4.64"
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("D17").Value = '34.864.95'
$ws.Range("E17").Value = '  -0.58%  '
$ws.Range("D18").Value = "'69.83"
$ws.Range("E18").Value = '  +0.14%  '
$ws.Range("E19").Value = '  -0.28%  '
$ws.Range("D20").Value = "'240.57"
$ws.Range("E20").Value = '  +1.06%  '
$ws.Range("D21").Value = "'12.17"
$ws.Range("E21").Value = '  +2.24%  '
$ws.Range("E22").Value = '  -0.48%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("E24").Value = '  +0.98%  '
$ws.Range("D25").Value = "'171.56"
$ws.Range("E25").Value = '  -0.32%  '
$ws.Range("E26").Value = '  -0.57%  '
$ws.Range("D27").Value = "'17.44"
$ws.Range("E27").Value = '  -0.34%  '
$ws.Range("E28").Value = '  +2.19%  '
$ws.Range("E29").Value = '  -4.19%  '
$ws.Range("E31").Value = '  -0.55%  '
$ws.Range("E32").Value = '  -4.38%  '
$ws.Range("E33").Value = '  -1.86%  '
$ws.Range("E34").Value = '  +6.72%  '
$ws.Range("E35").Value = '  +6.61%  '
$ws.Range("E36").Value = '  +12.85%  '
$ws.Range("D37").Value = "'0.693"
$ws.Range("E37").Value = '  +1.89%  '
$ws.Range("E38").Value = '  +7.04%  '
$ws.Range("D39").Value = "'90.41"
$ws.Range("E39").Value = '  -1.73%  '
$ws.Range("D40").Value = '1.345.43'
$ws.Range("E40").Value = '  +2.53%  '
$ws.Range("E41").Value = '  -0.17%  '
$ws.Range("D42").Value = "'14.91"
$ws.Range("E42").Value = '  +2.68%  '
$ws.Range("E43").Value = '  +0.59%  '
$ws.Range("E44").Value = '  -2.77%  '
$ws.Range("D45").Value = "'2.76"
$ws.Range("E45").Value = '  -0.27%  '
$ws.Range("D46").Value = "'6.26"
$ws.Range("E46").Value = '  -0.84%  '
$ws.Range("E47").Value = '  +2.13%  '
$ws.Range("D48").Value = '2.019.11'
$ws.Range("E48").Value = '  +1.45%  '
$ws.Range("E49").Value = '  +21.89%  '
$ws.Range("E51").Value = '  +1.56%  '

# Clear the quote-prefix formatting artifact introduced by forcing text values
foreach ($ref in @("D5", "D8", "D10", "D11", "D13", "D15", "D16", "D18", "D20", "D21", "D25", "D27", "D37", "D39", "D42", "D45", "D46")) {
    $ws.Range($ref).ClearFormats()
}
